# Update "想去人数" (interest count) figures in column F across all sheets,
# matching a refreshed data export (commit: "Update gh-pages to output
# generated at 456a3b4"). Only numeric values in column F change; nothing
# else in the workbook is modified.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 129
$ws.Range("F3").Value = 1312
$ws.Range("F4").Value = 1026
$ws.Range("F5").Value = 966
$ws.Range("F7").Value = 110
$ws.Range("F8").Value = 5
$ws.Range("F10").Value = 486
$ws.Range("F11").Value = 653
$ws.Range("F12").Value = 5
$ws.Range("F13").Value = 1860
$ws.Range("F14").Value = 4521
$ws.Range("F15").Value = 1309
$ws.Range("F17").Value = 2826
$ws.Range("F19").Value = 15
$ws.Range("F20").Value = 1139
$ws.Range("F21").Value = 3833
$ws.Range("F22").Value = 850
$ws.Range("F23").Value = 853
$ws.Range("F24").Value = 1539
$ws.Range("F25").Value = 52
$ws.Range("F26").Value = 2509
$ws.Range("F27").Value = 4
$ws.Range("F28").Value = 17
$ws.Range("F29").Value = 130
$ws.Range("F30").Value = 904
$ws.Range("F31").Value = 67
$ws.Range("F32").Value = 184
$ws.Range("F34").Value = 269
$ws.Range("F35").Value = 45
$ws.Range("F36").Value = 89
$ws.Range("F37").Value = 1467
$ws.Range("F38").Value = 2030
$ws.Range("F39").Value = 963
$ws.Range("F41").Value = 14
$ws.Range("F43").Value = 131
$ws.Range("F44").Value = 622
$ws.Range("F45").Value = 321
$ws.Range("F46").Value = 154
$ws.Range("F47").Value = 177
$ws.Range("F48").Value = 251
$ws.Range("F49").Value = 91

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 157
$ws.Range("F3").Value = 12
$ws.Range("F6").Value = 11
$ws.Range("F19").Value = 7

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 554

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 554
$ws.Range("F3").Value = 1312
$ws.Range("F4").Value = 1026
$ws.Range("F5").Value = 966
$ws.Range("F6").Value = 157
$ws.Range("F7").Value = 110
$ws.Range("F9").Value = 486
$ws.Range("F10").Value = 1860
$ws.Range("F11").Value = 4522
$ws.Range("F12").Value = 1309
$ws.Range("F14").Value = 11
$ws.Range("F15").Value = 2826
$ws.Range("F16").Value = 15
$ws.Range("F17").Value = 1139
$ws.Range("F18").Value = 3833
$ws.Range("F19").Value = 850
$ws.Range("F20").Value = 853
$ws.Range("F21").Value = 1539
$ws.Range("F23").Value = 52
$ws.Range("F24").Value = 2509
$ws.Range("F28").Value = 130
$ws.Range("F30").Value = 904
$ws.Range("F31").Value = 184
$ws.Range("F34").Value = 269
$ws.Range("F35").Value = 1467
$ws.Range("F36").Value = 2030
$ws.Range("F37").Value = 963
$ws.Range("F40").Value = 14
$ws.Range("F41").Value = 7
$ws.Range("F43").Value = 131
$ws.Range("F44").Value = 622
$ws.Range("F45").Value = 321
$ws.Range("F46").Value = 154
$ws.Range("F47").Value = 177
$ws.Range("F48").Value = 251
$ws.Range("F49").Value = 91
